# Apply updated cryptocurrency price/volume data to sheet1 (ActiveSheet).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.136.03"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +2.99%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.635.77"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +2.13%  "

# Row 4
$ws.Range("E4").Value = "  +0.08%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "595.92"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.00%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "155.99"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +4.01%  "

# Row 7
$ws.Range("E7").Value = "  +0.02%  "

# Row 8
$ws.Range("E8").Value = "  +1.21%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.117"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +7.31%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.400"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +4.58%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.78"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +1.47%  "

# Row 12
$ws.Range("E12").Value = "  +1.97%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "28.96"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +5.71%  "

# Row 14
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.106.53"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +2.07%  "

# Row 15
$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000181"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +16.85%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.041.33"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +3.21%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.680.23"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +3.90%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.53"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +2.48%  "

# Row 19
$ws.Range("E19").Value = "  +1.38%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "353.76"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +2.80%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.26"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +5.92%  "

# Row 22
$ws.Range("E22").Value = "  +0.15%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "68.18"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.34%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.70"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.01%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.50"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +2.79%  "

# Row 26
$ws.Range("E26").Value = "  -0.92%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.13"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +1.78%  "

# Row 28
$ws.Range("E28").Value = "  +1.21%  "

# Row 29
$ws.Range("B29").Value = "Binance-PegBSC-USD"
$ws.Range("C29").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.07%  "

# Row 30
$ws.Range("B30").Value = "PEPE"
$ws.Range("C30").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0942"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +11.65%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "524.86"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -6.94%  "

# Row 32
$ws.Range("E32").Value = "  +3.49%  "

# Row 33
$ws.Range("E33").Value = "  +1.95%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.65"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +8.27%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.32"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +3.52%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.427"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +3.79%  "

# Row 37
$ws.Range("B37").Value = "Monero"
$ws.Range("C37").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "164.94"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.91%  "

# Row 38
$ws.Range("B38").Value = "Stacks"
$ws.Range("C38").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.03"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +6.10%  "

# Row 39
$ws.Range("E39").Value = "  +3.71%  "

# Row 40
$ws.Range("E40").Value = "  +0.07%  "

# Row 41
$ws.Range("E41").Value = "  +0.10%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "42.21"
$ws.Range("D42").ClearFormats()

# Row 43
$ws.Range("E43").Value = "  -1.06%  "

# Row 44
$ws.Range("E44").Value = "  +3.92%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0616"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +5.79%  "

# Row 46
$ws.Range("E46").Value = "  +1.54%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.21"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +7.51%  "

# Row 48
$ws.Range("E48").Value = "  +3.11%  "

# Row 49
$ws.Range("E49").Value = "  +1.69%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0983"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +2.32%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.37"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +2.07%  "
